$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column retains text formatting so values like "0.3820" or
# "0.000007144" are not reinterpreted as numbers/scientific notation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.735.81'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.725.16'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").Value = '0.9977'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '240.25'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.4816'
$ws.Range("E7").Value = '  -1.44%  '
$ws.Range("D8").Value = '0.2588'
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").Value = '0.06182'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '1.722.51'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = '15.82'
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("D12").Value = '0.06855'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("D13").Value = '0.6028'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '4.457'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '76.82'
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = '0.9984'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '26.557.77'
$ws.Range("D18").Value = '0.9981'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '0.000007144'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '1.944.55'
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '4.410'
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '5.052'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = '139.70'
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +2.60%  '
$ws.Range("D28").Value = '106.26'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("D30").Value = '3.995'
$ws.Range("E30").Value = '  +2.30%  '
$ws.Range("D31").Value = '0.07912'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").Value = '3.663'
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").Value = '0.04529'
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("D34").Value = '2.594'
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").Value = '0.6167'
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").Value = '0.9294'
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("D38").Value = '2.451'
$ws.Range("E38").Value = '  +2.89%  '
$ws.Range("D39").Value = '1.990'
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("D40").Value = '0.9978'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("D42").Value = '5.599'
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("D43").Value = '99.77'
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").Value = '0.3820'
$ws.Range("E44").Value = '  -0.57%  '
$ws.Range("D45").Value = '6.767'
$ws.Range("E45").Value = '  -1.69%  '
$ws.Range("D46").Value = '0.1153'
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").Value = '0.05352'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '7.917'
$ws.Range("E48").Value = '  +3.29%  '
$ws.Range("D49").Value = '30.08'
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("D50").Value = '1.246'
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").Value = '51.43'
$ws.Range("E51").Value = '  +0.82%  '